$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 425
$ws.Range("I80").Value = 249.64285
$ws.Range("J80").Value = 775.7143
$ws.Range("K80").Value = 748.9285500000001
$ws.Range("L80").Value = 2327.1429
$ws.Range("M80").Value = 249.0714499999999
$ws.Range("N80").Value = -4323.1429

$ws.Range("H83").Value = 425
$ws.Range("I83").Value = 249.64285
$ws.Range("J83").Value = 775.7143
$ws.Range("K83").Value = 2246.78565
$ws.Range("L83").Value = 6981.428699999999
$ws.Range("M83").Value = 2745.21435
$ws.Range("N83").Value = -16965.4287

$ws.Range("H113").Value = 1649.5834
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1649.5834
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1649.5834
$ws.Range("N113").Value = -8157.5834
$ws.Range("M113").ClearContents()

$ws.Range("H138").Value = 2177227.5
$ws.Range("I138").Value = 4002653
$ws.Range("J138").Value = 4102.143
$ws.Range("K138").Value = 12007959
$ws.Range("L138").Value = 12306.429
$ws.Range("M138").Value = -12002819
$ws.Range("N138").Value = -22586.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 5186.8984
$ws.Range("I32").Value = 5555.851
$ws.Range("K32").Value = 5555.851
$ws.Range("M32").Value = -5268.851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2528.5881
$ws.Range("I86").Value = 2098.923
$ws.Range("J86").Value = 3925
$ws.Range("K86").Value = 2098.923
$ws.Range("L86").Value = 3925
$ws.Range("M86").Value = -975.9229999999998
$ws.Range("N86").Value = -6171

$ws.Range("H89").Value = 2528.5881
$ws.Range("I89").Value = 2098.923
$ws.Range("J89").Value = 3925
$ws.Range("K89").Value = 10494.615
$ws.Range("L89").Value = 19625
$ws.Range("M89").Value = -4878.614999999998
$ws.Range("N89").Value = -30857

$ws.Range("H105").Value = 4077.5
$ws.Range("I105").Value = 6155
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 6155
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -4408
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1017.8571
$ws.Range("I122").Value = 955.8333
$ws.Range("J122").Value = 1390
$ws.Range("K122").Value = 2867.4999
$ws.Range("L122").Value = 4170
$ws.Range("M122").Value = -417.4998999999998
$ws.Range("N122").Value = -9070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 10348050
$ws.Range("I80").Value = 24143188
$ws.Range("J80").Value = 1695
$ws.Range("K80").Value = 72429564
$ws.Range("L80").Value = 5085
$ws.Range("M80").Value = -72428628
$ws.Range("N80").Value = -6957

$ws.Range("H83").Value = 10348050
$ws.Range("I83").Value = 24143188
$ws.Range("J83").Value = 1695
$ws.Range("K83").Value = 217288692
$ws.Range("L83").Value = 15255
$ws.Range("M83").Value = -217284012
$ws.Range("N83").Value = -24615

$ws.Range("H92").Value = 2112.6
$ws.Range("I92").Value = 854.3333
$ws.Range("K92").Value = 2562.9999
$ws.Range("M92").Value = -1314.9999

$ws.Range("H113").Value = 662.3684
$ws.Range("I113").Value = 672.5
$ws.Range("J113").Value = 659.6667
$ws.Range("K113").Value = 2017.5
$ws.Range("L113").Value = 1979.0001
$ws.Range("M113").Value = 152.5
$ws.Range("N113").Value = -6319.0001

$ws.Range("H131").Value = 3988.7812
$ws.Range("J131").Value = 3111.6667
$ws.Range("L131").Value = 9335.000100000001
$ws.Range("N131").Value = -19415.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1286.2963
$ws.Range("I102").Value = 1162.8889
$ws.Range("J102").Value = 1533.1111
$ws.Range("K102").Value = 1162.8889
$ws.Range("L102").Value = 1533.1111
$ws.Range("M102").Value = 459.1111000000001
$ws.Range("N102").Value = -4777.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2598.375
$ws.Range("I7").Value = 2548
$ws.Range("J7").Value = 2648.75
$ws.Range("K7").Value = 2548
$ws.Range("L7").Value = 2648.75
$ws.Range("M7").Value = -2436
$ws.Range("N7").Value = -2872.75

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H61").Value = 1549.2727
$ws.Range("I61").Value = 1377.4286
$ws.Range("J61").Value = 1850
$ws.Range("K61").Value = 1377.4286
$ws.Range("L61").Value = 1850
$ws.Range("M61").Value = -1175.4286
$ws.Range("N61").Value = -2254

$ws.Range("H113").Value = 1549.2727
$ws.Range("I113").Value = 1377.4286
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 1377.4286
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 792.5714
$ws.Range("N113").Value = -6190

$ws.Range("H126").Value = 2598.375
$ws.Range("I126").Value = 2548
$ws.Range("J126").Value = 2648.75
$ws.Range("K126").Value = 7644
$ws.Range("L126").Value = 7946.25
$ws.Range("M126").Value = -5174
$ws.Range("N126").Value = -12886.25

$ws.Range("H136").Value = 2332.3572
$ws.Range("I136").Value = 1426.05
$ws.Range("J136").Value = 4598.125
$ws.Range("K136").Value = 4278.15
$ws.Range("L136").Value = 13794.375
$ws.Range("M136").Value = -1728.15
$ws.Range("N136").Value = -18894.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 413.89474
$ws.Range("I107").Value = 191.5
$ws.Range("K107").Value = 574.5
$ws.Range("M107").Value = 1345.5

$ws.Range("H126").Value = 4175.5557
$ws.Range("I126").Value = 4650
$ws.Range("J126").Value = 380
$ws.Range("K126").Value = 13950
$ws.Range("L126").Value = 1140
$ws.Range("M126").Value = -11480
$ws.Range("N126").Value = -6080

$ws.Range("H132").Value = 2097.4707
$ws.Range("I132").Value = 1923.2174
$ws.Range("J132").Value = 3700.6
$ws.Range("K132").Value = 5769.6522
$ws.Range("L132").Value = 11101.8
$ws.Range("M132").Value = -3239.6522
$ws.Range("N132").Value = -16161.8

$ws.Range("H140").Value = 21547.5
$ws.Range("J140").Value = 21547.5
$ws.Range("L140").Value = 21547.5
$ws.Range("N140").Value = -31907.5

Write-Host "Applied scheduled Sheets update to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR."
